# Resort the worksheet tabs: move "总计" (summary) ahead of "2022-Q2" (detail),
# so the sheet order becomes [总计, 2022-Q2] instead of [2022-Q2, 总计].
$wb = $excel.ActiveWorkbook

$summarySheet = $wb.Worksheets.Item("总计")
$detailSheet  = $wb.Worksheets.Item("2022-Q2")

# Move "总计" to be before "2022-Q2" (i.e., before the first sheet).
$summarySheet.Move($detailSheet)
